# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# All "Price"/"Volume(1h)" cells are stored as text (openpyxl-authored sheet),
# so for any new Price value that Excel would otherwise auto-parse as a
# number (e.g. "580.27", "1.00", "0.0000168"), we briefly force the cell to
# Text format before writing it and then ClearFormats() to drop the
# temporary number-format override again, leaving the cell's style untouched
# while keeping its value as literal text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.859.41"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "3.263.36"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.55"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "3.261.45"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.130"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.58"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.411"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").Value = "3.831.40"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.48"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.10%  "
$ws.Range("D16").Value = "67.877.03"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000168"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").Value = "3.304.74"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.61"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "393.09"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.514"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.188"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.52"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.95"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.68"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.51"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.99"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.26"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.26%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.12"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -4.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.90"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.79"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.53"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.45"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0687"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.44"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -8.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.62"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("D46").Value = "2.616.65"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.97"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "332.45"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.95%  "
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.101"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.24%  "
